# feat: add 2022-Q1 data
#
# 1. Duplicate the "总计" (totals) sheet; the duplicate keeps the "总计"
#    name/role (and gets the running data-update described in step 3),
#    while the original physical sheet is repurposed (renamed) into the new
#    "2022-Q1" quarterly snapshot sheet -- this keeps sheetPr/pageMargins
#    etc. consistent with the rest of the workbook's quarterly sheets.
# 2. Rebuild the new "2022-Q1" sheet with the same header layout as the
#    other quarterly sheets, plus the new fund snapshot row.
# 3. Insert a new top data row into the "总计" sheet for "2022-Q1", pushing
#    the existing quarters down one row and renumbering the index column.

$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count
$origTotal = $wb.Worksheets.Item($sheetCount)          # "总计" is the last sheet
$totalName = $origTotal.Name
$origIndex = $origTotal.Index
$templateSheet = $wb.Worksheets.Item($sheetCount - 1)  # most recent quarter sheet ("2021-Q4")

# --- 1. Duplicate "总计"; the copy becomes the (updated) totals sheet ---

$origTotal.Copy($null, $origTotal)
$totalSheet = $wb.Worksheets.Item($origIndex + 1)

$newSheet = $origTotal                 # repurpose the original sheet object
$newSheet.Name = "2022-Q1"
$totalSheet.Name = $totalName

# --- 2. Rebuild the new "2022-Q1" quarterly sheet -----------------------

$newSheet.UsedRange.Delete()

# Reuse the header formatting/style from the template quarter sheet, then
# (re)write all the header text explicitly.
$templateSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Reuse the row-2 index-cell style from the template quarter sheet.
$templateSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$newSheet.Range("A2").Value = 0

# B2/D2/E2/F2/G2 hold numeric-looking figures that must stay TEXT (matching
# the other quarter sheets), so force a text format while writing them, then
# drop the number-format override again (keeps the default, unstyled cell).
$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("B2").Value = "513080"
$newSheet.Range("C2").Value = "华安法国CAC40ETF（QDII）"
$newSheet.Range("D2").Value = "0.60"
$newSheet.Range("E2").Value = "96.69"
$newSheet.Range("F2").Value = "3.56"
$newSheet.Range("G2").Value = "0.0214"
$newSheet.Range("B2:G2").ClearFormats()
$newSheet.Range("H2").Value = 8

# --- 3. Insert the new top row into the "总计" sheet ---------------------

$lastRow = $totalSheet.Cells.Item($totalSheet.Rows.Count, 1).End(-4162).Row

$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A2:D2").ClearFormats()

# Reuse the index-cell style (column A) from the row that just got pushed
# down to row 3.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.02

# Renumber the index column (A) for the rows that shifted down.
for ($r = 3; $r -le ($lastRow + 1); $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}
